# refatoração - cálculos de apoio médio
#
# Splits the single "apoio_medio" aggregate block and the
# "contribuicoes"/"media_contribuicoes" block into fuller avg/std/min/max
# stat sets (mirroring the existing arrecadado_* columns), and renames the
# old media/std/min/max_sucesso headers to arrecadado_avg/std/min/max.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns after "apoio_medio" (L) for apoio_std/min/max,
# and 3 new columns after "media_contribuicoes" (now Q, contribuicoes_med)
# for contribuicoes_std/min/max. Inserting picks up the left neighbour's
# column style automatically (style 3 for M:O, style 1 for R:T), matching
# the surrounding arrecadado_* / contribuicoes columns.
$ws.Columns("M:O").Insert()
$ws.Columns("R:T").Insert()

# Rename the arrecadado_sucesso stat columns.
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"

# New apoio_medio stat columns.
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"

# Rename + extend the contribuicoes stat columns.
$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"

# Updated / new data values for row 2.
$ws.Range("L2").Value = 91.85574933975617
$ws.Range("M2").Value = 49.08980856017526
$ws.Range("N2").Value = 13.93896149503088
$ws.Range("O2").Value = 792.0360759681182
$ws.Range("R2").Value = 423.019225146675
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 6494
